{"js": "// The author applied direct character formatting (font size 12pt) to the\n// substring \"LIFO (Last-in, First-Out), \" inside the paragraph that begins\n// \"Strack l\u00e0 m\u1ed9t c\u1ea5u tr\u00fac d\u1eef li\u1ec7u tuy\u1ebfn t\u00ednh...\". This splits the original\n// run(s) so the resized text becomes its own run(s).\nconst body = context.document.body;\n\nconst results = body.search(\"LIFO (Last-in, First-Out), \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target text \"LIFO (Last-in, First-Out), \" not found.');\n}\n\nconst target = results.items[0];\ntarget.font.size = 12;\nawait context.sync();\n", "ps1": "# The author applied direct character formatting (font size 12pt) to the\n# substring \"LIFO (Last-in, First-Out), \" inside the paragraph that begins\n# \"Strack l\u00e0 m\u1ed9t c\u1ea5u tr\u00fac d\u1eef li\u1ec7u tuy\u1ebfn t\u00ednh...\". This splits the original\n# run(s) so the resized text becomes its own run(s).\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"LIFO (Last-in, First-Out), \"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found -and $find.Found) {\n    $rng.Font.Size = 12\n}\n"}
